$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-40 with shuffled/updated data per diff
$ws.Range('D2').Value = 44329
$ws.Range('L2').Value = 'Primera'
$ws.Range('M2').Value = 50
$ws.Range('N2').Value = 16000
$ws.Range('O2').Value = 16000
$ws.Range('P2').Value = 16000
$ws.Range('Q2').Value = '$/bandeja 10 kilos'
$ws.Range('S2').Value = 1600
$ws.Range('T2').Value = 10

$ws.Range('D3').Value = 44396
$ws.Range('L3').Value = 'Primera'
$ws.Range('M3').Value = 60
$ws.Range('N3').Value = 17000
$ws.Range('O3').Value = 17000
$ws.Range('P3').Value = 17000
$ws.Range('Q3').Value = '$/bandeja 10 kilos'
$ws.Range('S3').Value = 1700
$ws.Range('T3').Value = 10

$ws.Range('D4').Value = 44396
$ws.Range('L4').Value = 'Segunda'
$ws.Range('M4').Value = 56
$ws.Range('N4').Value = 15000
$ws.Range('O4').Value = 15000
$ws.Range('P4').Value = 15000
$ws.Range('Q4').Value = '$/bandeja 10 kilos'
$ws.Range('S4').Value = 1500
$ws.Range('T4').Value = 10

$ws.Range('D5').Value = 44354
$ws.Range('L5').Value = 'Primera'
$ws.Range('M5').Value = 45
$ws.Range('N5').Value = 15000
$ws.Range('O5').Value = 15000
$ws.Range('P5').Value = 15000
$ws.Range('Q5').Value = '$/bandeja 10 kilos'
$ws.Range('S5').Value = 1500
$ws.Range('T5').Value = 10

$ws.Range('D6').Value = 44398
$ws.Range('L6').Value = 'Primera'
$ws.Range('M6').Value = 60
$ws.Range('N6').Value = 17000
$ws.Range('O6').Value = 17000
$ws.Range('P6').Value = 17000
$ws.Range('Q6').Value = '$/bandeja 10 kilos'
$ws.Range('S6').Value = 1700
$ws.Range('T6').Value = 10

$ws.Range('D7').Value = 44398
$ws.Range('L7').Value = 'Segunda'
$ws.Range('M7').Value = 50
$ws.Range('N7').Value = 15000
$ws.Range('O7').Value = 15000
$ws.Range('P7').Value = 15000
$ws.Range('Q7').Value = '$/bandeja 10 kilos'
$ws.Range('S7').Value = 1500
$ws.Range('T7').Value = 10

$ws.Range('D8').Value = 44441
$ws.Range('L8').Value = 'Primera'
$ws.Range('M8').Value = 80
$ws.Range('N8').Value = 15000
$ws.Range('O8').Value = 15000
$ws.Range('P8').Value = 15000
$ws.Range('Q8').Value = '$/bandeja 10 kilos'
$ws.Range('S8').Value = 1500
$ws.Range('T8').Value = 10

$ws.Range('D9').Value = 44420
$ws.Range('L9').Value = 'Primera'
$ws.Range('M9').Value = 54
$ws.Range('N9').Value = 18000
$ws.Range('O9').Value = 18000
$ws.Range('P9').Value = 18000
$ws.Range('Q9').Value = '$/bandeja 10 kilos'
$ws.Range('S9').Value = 1800
$ws.Range('T9').Value = 10

$ws.Range('D10').Value = 44420
$ws.Range('L10').Value = 'Segunda'
$ws.Range('M10').Value = 50
$ws.Range('N10').Value = 15000
$ws.Range('O10').Value = 15000
$ws.Range('P10').Value = 15000
$ws.Range('Q10').Value = '$/bandeja 10 kilos'
$ws.Range('S10').Value = 1500
$ws.Range('T10').Value = 10

$ws.Range('D11').Value = 44323
$ws.Range('L11').Value = 'Primera'
$ws.Range('M11').Value = 48
$ws.Range('N11').Value = 24000
$ws.Range('O11').Value = 24000
$ws.Range('P11').Value = 24000
$ws.Range('Q11').Value = '$/caja 15 kilos granel'
$ws.Range('S11').Value = 1600
$ws.Range('T11').Value = 15

$ws.Range('D12').Value = 44431
$ws.Range('L12').Value = 'Primera'
$ws.Range('M12').Value = 65
$ws.Range('N12').Value = 18000
$ws.Range('O12').Value = 18000
$ws.Range('P12').Value = 18000
$ws.Range('Q12').Value = '$/bandeja 10 kilos'
$ws.Range('S12').Value = 1800
$ws.Range('T12').Value = 10

$ws.Range('D13').Value = 44431
$ws.Range('L13').Value = 'Segunda'
$ws.Range('M13').Value = 60
$ws.Range('N13').Value = 16000
$ws.Range('O13').Value = 16000
$ws.Range('P13').Value = 16000
$ws.Range('Q13').Value = '$/bandeja 10 kilos'
$ws.Range('S13').Value = 1600
$ws.Range('T13').Value = 10

$ws.Range('D14').Value = 44371
$ws.Range('L14').Value = 'Primera'
$ws.Range('M14').Value = 67
$ws.Range('N14').Value = 18000
$ws.Range('O14').Value = 18000
$ws.Range('P14').Value = 18000
$ws.Range('Q14').Value = '$/bandeja 10 kilos'
$ws.Range('S14').Value = 1800
$ws.Range('T14').Value = 10

$ws.Range('D15').Value = 44417
$ws.Range('L15').Value = 'Primera'
$ws.Range('M15').Value = 56
$ws.Range('N15').Value = 16000
$ws.Range('O15').Value = 16000
$ws.Range('P15').Value = 16000
$ws.Range('Q15').Value = '$/bandeja 10 kilos'
$ws.Range('S15').Value = 1600
$ws.Range('T15').Value = 10

$ws.Range('D16').Value = 44417
$ws.Range('L16').Value = 'Segunda'
$ws.Range('M16').Value = 60
$ws.Range('N16').Value = 14000
$ws.Range('O16').Value = 14000
$ws.Range('P16').Value = 14000
$ws.Range('Q16').Value = '$/bandeja 10 kilos'
$ws.Range('S16').Value = 1400
$ws.Range('T16').Value = 10

$ws.Range('D17').Value = 44382
$ws.Range('L17').Value = 'Primera'
$ws.Range('M17').Value = 58
$ws.Range('N17').Value = 17000
$ws.Range('O17').Value = 17000
$ws.Range('P17').Value = 17000
$ws.Range('Q17').Value = '$/bandeja 10 kilos'
$ws.Range('S17').Value = 1700
$ws.Range('T17').Value = 10

$ws.Range('D18').Value = 44370
$ws.Range('L18').Value = 'Primera'
$ws.Range('M18').Value = 50
$ws.Range('N18').Value = 17000
$ws.Range('O18').Value = 17000
$ws.Range('P18').Value = 17000
$ws.Range('Q18').Value = '$/bandeja 10 kilos'
$ws.Range('S18').Value = 1700
$ws.Range('T18').Value = 10

$ws.Range('D19').Value = 44315
$ws.Range('L19').Value = 'Primera'
$ws.Range('M19').Value = 60
$ws.Range('N19').Value = 24000
$ws.Range('O19').Value = 24000
$ws.Range('P19').Value = 24000
$ws.Range('Q19').Value = '$/caja 15 kilos granel'
$ws.Range('S19').Value = 1600
$ws.Range('T19').Value = 15

$ws.Range('D20').Value = 44454
$ws.Range('L20').Value = 'Primera'
$ws.Range('M20').Value = 45
$ws.Range('N20').Value = 15000
$ws.Range('O20').Value = 15000
$ws.Range('P20').Value = 15000
$ws.Range('Q20').Value = '$/bandeja 10 kilos'
$ws.Range('S20').Value = 1500
$ws.Range('T20').Value = 10

$ws.Range('D21').Value = 44410
$ws.Range('L21').Value = 'Primera'
$ws.Range('M21').Value = 75
$ws.Range('N21').Value = 15000
$ws.Range('O21').Value = 15000
$ws.Range('P21').Value = 15000
$ws.Range('Q21').Value = '$/bandeja 10 kilos'
$ws.Range('S21').Value = 1500
$ws.Range('T21').Value = 10

$ws.Range('D22').Value = 44319
$ws.Range('L22').Value = 'Primera'
$ws.Range('M22').Value = 60
$ws.Range('N22').Value = 24000
$ws.Range('O22').Value = 24000
$ws.Range('P22').Value = 24000
$ws.Range('Q22').Value = '$/caja 15 kilos granel'
$ws.Range('S22').Value = 1600
$ws.Range('T22').Value = 15

$ws.Range('D23').Value = 44473
$ws.Range('L23').Value = 'Primera'
$ws.Range('M23').Value = 85
$ws.Range('N23').Value = 18000
$ws.Range('O23').Value = 18000
$ws.Range('P23').Value = 18000
$ws.Range('Q23').Value = '$/bandeja 10 kilos'
$ws.Range('S23').Value = 1800
$ws.Range('T23').Value = 10

$ws.Range('D24').Value = 44413
$ws.Range('L24').Value = 'Primera'
$ws.Range('M24').Value = 60
$ws.Range('N24').Value = 15000
$ws.Range('O24').Value = 15000
$ws.Range('P24').Value = 15000
$ws.Range('Q24').Value = '$/bandeja 10 kilos'
$ws.Range('S24').Value = 1500
$ws.Range('T24').Value = 10

$ws.Range('D25').Value = 44413
$ws.Range('L25').Value = 'Segunda'
$ws.Range('M25').Value = 58
$ws.Range('N25').Value = 13000
$ws.Range('O25').Value = 13000
$ws.Range('P25').Value = 13000
$ws.Range('Q25').Value = '$/bandeja 10 kilos'
$ws.Range('S25').Value = 1300
$ws.Range('T25').Value = 10

$ws.Range('D26').Value = 44469
$ws.Range('L26').Value = 'Primera'
$ws.Range('M26').Value = 50
$ws.Range('N26').Value = 16000
$ws.Range('O26').Value = 16000
$ws.Range('P26').Value = 16000
$ws.Range('Q26').Value = '$/bandeja 10 kilos'
$ws.Range('S26').Value = 1600
$ws.Range('T26').Value = 10

$ws.Range('D27').Value = 44445
$ws.Range('L27').Value = 'Primera'
$ws.Range('M27').Value = 68
$ws.Range('N27').Value = 15000
$ws.Range('O27').Value = 15000
$ws.Range('P27').Value = 15000
$ws.Range('Q27').Value = '$/bandeja 10 kilos'
$ws.Range('S27').Value = 1500
$ws.Range('T27').Value = 10

$ws.Range('D28').Value = 44391
$ws.Range('L28').Value = 'Primera'
$ws.Range('M28').Value = 50
$ws.Range('N28').Value = 17000
$ws.Range('O28').Value = 17000
$ws.Range('P28').Value = 17000
$ws.Range('Q28').Value = '$/bandeja 10 kilos'
$ws.Range('S28').Value = 1700
$ws.Range('T28').Value = 10

$ws.Range('D29').Value = 44391
$ws.Range('L29').Value = 'Segunda'
$ws.Range('M29').Value = 45
$ws.Range('N29').Value = 15000
$ws.Range('O29').Value = 15000
$ws.Range('P29').Value = 15000
$ws.Range('Q29').Value = '$/bandeja 10 kilos'
$ws.Range('S29').Value = 1500
$ws.Range('T29').Value = 10

$ws.Range('D30').Value = 44475
$ws.Range('L30').Value = 'Primera'
$ws.Range('M30').Value = 56
$ws.Range('N30').Value = 17000
$ws.Range('O30').Value = 17000
$ws.Range('P30').Value = 17000
$ws.Range('Q30').Value = '$/bandeja 10 kilos'
$ws.Range('S30').Value = 1700
$ws.Range('T30').Value = 10

$ws.Range('D31').Value = 44385
$ws.Range('L31').Value = 'Primera'
$ws.Range('M31').Value = 60
$ws.Range('N31').Value = 17000
$ws.Range('O31').Value = 17000
$ws.Range('P31').Value = 17000
$ws.Range('Q31').Value = '$/bandeja 10 kilos'
$ws.Range('S31').Value = 1700
$ws.Range('T31').Value = 10

$ws.Range('D32').Value = 44385
$ws.Range('L32').Value = 'Segunda'
$ws.Range('M32').Value = 50
$ws.Range('N32').Value = 15000
$ws.Range('O32').Value = 15000
$ws.Range('P32').Value = 15000
$ws.Range('Q32').Value = '$/bandeja 10 kilos'
$ws.Range('S32').Value = 1500
$ws.Range('T32').Value = 10

$ws.Range('D33').Value = 44453
$ws.Range('L33').Value = 'Primera'
$ws.Range('M33').Value = 50
$ws.Range('N33').Value = 15000
$ws.Range('O33').Value = 15000
$ws.Range('P33').Value = 15000
$ws.Range('Q33').Value = '$/bandeja 10 kilos'
$ws.Range('S33').Value = 1500
$ws.Range('T33').Value = 10

$ws.Range('D34').Value = 44435
$ws.Range('L34').Value = 'Primera'
$ws.Range('M34').Value = 115
$ws.Range('N34').Value = 18000
$ws.Range('O34').Value = 18000
$ws.Range('P34').Value = 18000
$ws.Range('Q34').Value = '$/bandeja 10 kilos'
$ws.Range('S34').Value = 1800
$ws.Range('T34').Value = 10

$ws.Range('D35').Value = 44435
$ws.Range('L35').Value = 'Segunda'
$ws.Range('M35').Value = 60
$ws.Range('N35').Value = 16000
$ws.Range('O35').Value = 16000
$ws.Range('P35').Value = 16000
$ws.Range('Q35').Value = '$/bandeja 10 kilos'
$ws.Range('S35').Value = 1600
$ws.Range('T35').Value = 10

$ws.Range('D36').Value = 44433
$ws.Range('L36').Value = 'Primera'
$ws.Range('M36').Value = 50
$ws.Range('N36').Value = 18000
$ws.Range('O36').Value = 18000
$ws.Range('P36').Value = 18000
$ws.Range('Q36').Value = '$/bandeja 10 kilos'
$ws.Range('S36').Value = 1800
$ws.Range('T36').Value = 10

$ws.Range('D37').Value = 44468
$ws.Range('L37').Value = 'Primera'
$ws.Range('M37').Value = 50
$ws.Range('N37').Value = 16000
$ws.Range('O37').Value = 16000
$ws.Range('P37').Value = 16000
$ws.Range('Q37').Value = '$/bandeja 10 kilos'
$ws.Range('S37').Value = 1600
$ws.Range('T37').Value = 10

$ws.Range('D38').Value = 44321
$ws.Range('L38').Value = 'Primera'
$ws.Range('M38').Value = 42
$ws.Range('N38').Value = 24000
$ws.Range('O38').Value = 24000
$ws.Range('P38').Value = 24000
$ws.Range('Q38').Value = '$/caja 15 kilos granel'
$ws.Range('S38').Value = 1600
$ws.Range('T38').Value = 15

$ws.Range('D39').Value = 44466
$ws.Range('L39').Value = 'Primera'
$ws.Range('M39').Value = 70
$ws.Range('N39').Value = 16000
$ws.Range('O39').Value = 16000
$ws.Range('P39').Value = 16000
$ws.Range('Q39').Value = '$/bandeja 10 kilos'
$ws.Range('S39').Value = 1600
$ws.Range('T39').Value = 10

$ws.Range('D40').Value = 44389
$ws.Range('L40').Value = 'Primera'
$ws.Range('M40').Value = 60
$ws.Range('N40').Value = 17000
$ws.Range('O40').Value = 17000
$ws.Range('P40').Value = 17000
$ws.Range('Q40').Value = '$/bandeja 10 kilos'
$ws.Range('S40').Value = 1700
$ws.Range('T40').Value = 10

# Add new row 41
$ws.Range('A41').Value = 3
$ws.Range('B41').Value = 'Femacal de La Calera'
$ws.Range('C41').Value = 'Coquimbo'
$ws.Range('D41').Value = 44389
$ws.Range('D41').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range('E41').Value = 5
$ws.Range('F41').Value = 'Fruta'
$ws.Range('G41').Value = 100108
$ws.Range('H41').Value = 'Tropicales y subtropicales'
$ws.Range('I41').Value = 100108004
$ws.Range('J41').Value = 'Papaya'
$ws.Range('K41').Value = 'Cultivar IV Región'
$ws.Range('L41').Value = 'Segunda'
$ws.Range('M41').Value = 50
$ws.Range('N41').Value = 15000
$ws.Range('O41').Value = 15000
$ws.Range('P41').Value = 15000
$ws.Range('Q41').Value = '$/bandeja 10 kilos'
$ws.Range('R41').Value = 'Provincia del Elquí'
$ws.Range('S41').Value = 1500
$ws.Range('T41').Value = 10
